# Correct wrong scan numbers in info file.
# Rows 4-9 in column B ("scandir") incorrectly referenced scan numbers
# 065-070, but they should all reference scan number 064, matching
# rows 2-3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$correctValue = "064_GMB_ringscan_500mA_Cu100_gap5p7_orca_10x_dist50mm_100ms"

$ws.Range("B4").Value = $correctValue
$ws.Range("B5").Value = $correctValue
$ws.Range("B6").Value = $correctValue
$ws.Range("B7").Value = $correctValue
$ws.Range("B8").Value = $correctValue
$ws.Range("B9").Value = $correctValue

# Update the active cell/selection to B10, matching the saved view state.
$ws.Range("B10").Select()
